$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "27.864.61"
Set-TextValue "E2" "  +0.94%  "
Set-TextValue "D3" "1.886.87"
Set-TextValue "E3" "  +0.79%  "
Set-TextValue "E4" "  +1.67%  "
Set-TextValue "E5" "  +1.78%  "
Set-TextValue "E6" "  +1.65%  "
Set-TextValue "D7" "0.4685"
Set-TextValue "E7" "  -0.32%  "
Set-TextValue "D8" "0.3913"
Set-TextValue "E8" "  -1.51%  "
Set-TextValue "D9" "47.92"
Set-TextValue "E9" "  +4.47%  "
Set-TextValue "D10" "0.08046"
Set-TextValue "E10" "  -0.37%  "
Set-TextValue "D11" "1.018"
Set-TextValue "E11" "  -0.65%  "
Set-TextValue "D12" "21.86"
Set-TextValue "E12" "  +0.19%  "
Set-TextValue "D13" "1.868.89"
Set-TextValue "D14" "5.949"
Set-TextValue "E14" "  -0.07%  "
Set-TextValue "D15" "7.080"
Set-TextValue "E15" "  -1.61%  "
Set-TextValue "E16" "  +1.80%  "
Set-TextValue "D17" "0.06768"
Set-TextValue "E17" "  +2.99%  "
Set-TextValue "D18" "87.25"
Set-TextValue "E18" "  +0.33%  "
Set-TextValue "D19" "0.00001047"
Set-TextValue "E19" "  +0.44%  "
Set-TextValue "D20" "17.17"
Set-TextValue "E20" "  -0.57%  "
Set-TextValue "E21" "  +1.52%  "
Set-TextValue "D22" "27.914.25"
Set-TextValue "E22" "  +1.14%  "
Set-TextValue "D23" "5.494"
Set-TextValue "E23" "  -0.37%  "
Set-TextValue "E24" "  -0.23%  "
Set-TextValue "D25" "2.343"
Set-TextValue "E25" "  +1.77%  "
Set-TextValue "D26" "2.118.63"
Set-TextValue "E26" "  +2.10%  "
Set-TextValue "D27" "160.06"
Set-TextValue "E27" "  +3.63%  "
Set-TextValue "D28" "20.07"
Set-TextValue "E28" "  -1.06%  "
Set-TextValue "D29" "2.071"
Set-TextValue "E29" "  -0.81%  "
Set-TextValue "D30" "5.458"
Set-TextValue "E30" "  -1.70%  "
Set-TextValue "D31" "121.70"
Set-TextValue "E31" "  -0.63%  "
Set-TextValue "D32" "0.9647"
Set-TextValue "E32" "  +0.90%  "
Set-TextValue "D33" "0.09477"
Set-TextValue "E33" "  -0.12%  "
Set-TextValue "D34" "3.646"
Set-TextValue "E34" "  +1.26%  "
Set-TextValue "D35" "1.407"
Set-TextValue "E35" "  -4.75%  "
Set-TextValue "D36" "5.335"
Set-TextValue "E36" "  +0.49%  "
Set-TextValue "D37" "0.06109"
Set-TextValue "E37" "  +0.06%  "
Set-TextValue "D38" "0.02250"
Set-TextValue "E38" "  -0.26%  "
Set-TextValue "D39" "1.217"
Set-TextValue "E39" "  -0.22%  "
Set-TextValue "D40" "8.030"
Set-TextValue "E40" "  -1.60%  "
Set-TextValue "E41" "  -0.47%  "
Set-TextValue "D42" "0.1882"
Set-TextValue "E42" "  -1.03%  "
Set-TextValue "D43" "10.26"
Set-TextValue "E43" "  -0.83%  "
Set-TextValue "D44" "1.269"
Set-TextValue "E44" "  +1.27%  "
Set-TextValue "D45" "0.5690"
Set-TextValue "E45" "  -0.29%  "
Set-TextValue "D46" "12.19"
Set-TextValue "E46" "  -0.43%  "
Set-TextValue "D47" "3.405"
Set-TextValue "E47" "  -0.07%  "
Set-TextValue "D48" "1.924"
Set-TextValue "E48" "  -0.82%  "
Set-TextValue "D49" "0.06929"
Set-TextValue "E49" "  +1.86%  "
Set-TextValue "D50" "114.02"
Set-TextValue "E50" "  +3.84%  "
Set-TextValue "D51" "1.066"
Set-TextValue "E51" "  +0.28%  "
